$d = $word.ActiveDocument

# Word COM Find.Execute positional signature:
#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace
# Wrap: 0 = wdFindStop, 1 = wdFindContinue, 2 = wdFindAsk
# Replace: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll

# --- 1) "Login to MoDaC prod server (ssh fsdsgl-modac01p.ncifcrf.gov)" ---
# The ".ncifcrf.gov" run and the ")" run get merged into a single run
# with text ".ncifcrf.gov)". Scope the search to that single paragraph and
# only touch the ".ncifcrf.gov)" tail so the colored host-name run
# (fsdsgl-modac01p) keeps its own separate run/formatting.
$p1 = $d.Paragraphs.Item(151)
Write-Output ("p1: " + $p1.Range.Text)
$r1 = $p1.Range.Duplicate
$ok1 = $r1.Find.Execute(".ncifcrf.gov)", $true, $false, $false, $false, $false, $true, 0, $false, ".ncifcrf.gov)", 2)
Write-Output ("step1 replaced => " + $ok1)

# --- 2) "Go to the path: /usr/share/tomcat/webapps/web-doe-<version>/WEB-INF/classes" ---
# The three path-fragment runs ("/share/tomcat/webapps/web-doe-",
# "<version>", "/WEB-INF/classes") get merged into one run. Scope to the
# specific paragraph, leaving the earlier "/" and "usr" runs untouched.
$p2 = $d.Paragraphs.Item(153)
Write-Output ("p2: " + $p2.Range.Text)
$r2 = $p2.Range.Duplicate
$ok2 = $r2.Find.Execute("/share/tomcat/webapps/web-doe-<version>/WEB-INF/classes", $true, $false, $false, $false, $false, $true, 0, $false, "/share/tomcat/webapps/web-doe-<version>/WEB-INF/classes", 2)
Write-Output ("step2 replaced => " + $ok2)

# --- 3) "Save the file " ---
# Trailing space is removed from the run, and a new run containing "." is
# appended right after it (same color/size, but its own <w:r>).
$p3 = $d.Paragraphs.Item(156)
Write-Output ("p3: " + $p3.Range.Text)
$r3 = $p3.Range.Duplicate
$ok3 = $r3.Find.Execute("Save the file ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
Write-Output ("step3 found => " + $ok3)
if ($ok3) {
  $r3.Text = "Save the file"
  $r3.Collapse(0)
  $r3.InsertAfter(".")
  # Toggle a character property on/off so the engine keeps the new "."
  # text as its own run instead of silently re-merging it into the
  # preceding run just because the resulting formatting is identical.
  $r3.Font.Bold = 1
  $r3.Font.Bold = 0
}

Write-Output ("final p3: " + $d.Paragraphs.Item(156).Range.Text)
